$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 8.660232485948974, 1, 11.16088251738829)
    3 = @(1.459612070389937, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 4.429675500412797)
    4 = @(3.230985683306322, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1, 4.198080582305154)
    5 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 9.295990156953671)
    6 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    7 = @(1.459612070389937, 1.667794583268128, 26.21740644021617, 0.496779210170732, 0, 29.84159230404497)
    8 = @(0.3048080303191223, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 0, 1.649481363816475)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
